$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

# --- Update the StatQuery Cypher text (shared by CasesTab/SamplesTab/FilesTab rows) ---
$newQuery = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f:file)-[*]->(c)`nOPTIONAL MATCH (sf:file)-->(s)`nWITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p`nWHERE demo.breed IN ['American Staffordshire Terrier']`nRETURN  `n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct f) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"
$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# --- Update the active selection / scroll position on the sheet ---
$ws.Range("B2").Select()
